$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the data row (row 2, columns A:P) content while keeping formatting/styles
$ws.Range("A2:P2").ClearContents()

# Match the updated row heights from the target revision
$ws.Rows.Item(1).RowHeight = 91.5
$ws.Rows.Item(2).RowHeight = 15

# Update the view state (zoom + active selection) to match the target
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("A3").Select()
